$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value2 = 0.171830985915493
$ws.Range("C2").Value2 = 0.6084507042253521
$ws.Range("J2").Value2 = 0.01690140845070422
$ws.Range("P2").Value2 = 0.1295774647887324
$ws.Range("S2").Value2 = 0.07323943661971831
$ws.Range("B3").Value2 = 0.004329004329004329
$ws.Range("C3").Value2 = 0.03896103896103896
$ws.Range("J3").Value2 = 0.02597402597402598
$ws.Range("P3").Value2 = 0.7835497835497836
$ws.Range("S3").Value2 = 0.1471861471861472
$ws.Range("J4").Value2 = 0.05714285714285714
$ws.Range("P4").Value2 = 0.6857142857142857
$ws.Range("S4").Value2 = 0.2571428571428571
$ws.Range("J5").Value2 = 0.4
$ws.Range("P5").Value2 = 0.4
$ws.Range("S5").Value2 = 0.2
$ws.Range("B6").Value2 = 0.1008403361344538
$ws.Range("D6").Value2 = 0.02941176470588235
$ws.Range("E6").Value2 = 0.004201680672268907
$ws.Range("F6").Value2 = 0.09663865546218488
$ws.Range("J6").Value2 = 0.2352941176470588
$ws.Range("O6").Value2 = 0.02521008403361345
$ws.Range("Q6").Value2 = 0.1974789915966386
$ws.Range("R6").Value2 = 0.07142857142857142
$ws.Range("S6").Value2 = 0.2394957983193277
$ws.Range("B7").Value2 = 0.132295719844358
$ws.Range("D7").Value2 = 0.0311284046692607
$ws.Range("F7").Value2 = 0.06614785992217899
$ws.Range("J7").Value2 = 0.1439688715953307
$ws.Range("O7").Value2 = 0.02334630350194553
$ws.Range("Q7").Value2 = 0.2373540856031128
$ws.Range("R7").Value2 = 0.05836575875486381
$ws.Range("S7").Value2 = 0.3073929961089494
$ws.Range("B8").Value2 = 0.08663366336633663
$ws.Range("D8").Value2 = 0.03217821782178218
$ws.Range("E8").Value2 = 0.002475247524752475
$ws.Range("F8").Value2 = 0.0594059405940594
$ws.Range("J8").Value2 = 0.103960396039604
$ws.Range("O8").Value2 = 0.01237623762376238
$ws.Range("Q8").Value2 = 0.1633663366336634
$ws.Range("R8").Value2 = 0.1014851485148515
$ws.Range("S8").Value2 = 0.4381188118811881
$ws.Range("B9").Value2 = 0.0903954802259887
$ws.Range("D9").Value2 = 0.01129943502824859
$ws.Range("E9").Value2 = 0.005649717514124294
$ws.Range("F9").Value2 = 0.06779661016949153
$ws.Range("J9").Value2 = 0.1355932203389831
$ws.Range("O9").Value2 = 0.01694915254237288
$ws.Range("Q9").Value2 = 0.1751412429378531
$ws.Range("R9").Value2 = 0.096045197740113
$ws.Range("S9").Value2 = 0.4011299435028249
$ws.Range("B10").Value2 = 0.1298793470546487
$ws.Range("D10").Value2 = 0.028388928317956
$ws.Range("E10").Value2 = 0.0021291696238467
$ws.Range("F10").Value2 = 0.05748757984386089
$ws.Range("J10").Value2 = 0.1298793470546487
$ws.Range("O10").Value2 = 0.0134847409510291
$ws.Range("Q10").Value2 = 0.2178850248403123
$ws.Range("R10").Value2 = 0.09013484740951029
$ws.Range("S10").Value2 = 0.3307310149041874
$ws.Range("F11").Value2 = 0.002551020408163265
$ws.Range("G11").Value2 = 0.1428571428571428
$ws.Range("J11").Value2 = 0.07908163265306123
$ws.Range("K11").Value2 = 0.1862244897959184
$ws.Range("L11").Value2 = 0.576530612244898
$ws.Range("S11").Value2 = 0.01275510204081633
$ws.Range("G12").Value2 = 0.7725321888412017
$ws.Range("J12").Value2 = 0.1630901287553648
$ws.Range("K12").Value2 = 0.008583690987124463
$ws.Range("L12").Value2 = 0.02575107296137339
$ws.Range("S12").Value2 = 0.03004291845493562
$ws.Range("G13").Value2 = 0.68
$ws.Range("J13").Value2 = 0.3
$ws.Range("S13").Value2 = 0.02
$ws.Range("F15").Value2 = 0.04888888888888889
$ws.Range("H15").Value2 = 0.1511111111111111
$ws.Range("I15").Value2 = 0.06666666666666667
$ws.Range("J15").Value2 = 0.3555555555555556
$ws.Range("K15").Value2 = 0.06666666666666667
$ws.Range("M15").Value2 = 0.01333333333333333
$ws.Range("O15").Value2 = 0.12
$ws.Range("S15").Value2 = 0.1777777777777778
$ws.Range("F16").Value2 = 0.007547169811320755
$ws.Range("H16").Value2 = 0.1283018867924528
$ws.Range("I16").Value2 = 0.06415094339622641
$ws.Range("J16").Value2 = 0.4716981132075472
$ws.Range("K16").Value2 = 0.1471698113207547
$ws.Range("M16").Value2 = 0.02264150943396226
$ws.Range("O16").Value2 = 0.0339622641509434
$ws.Range("S16").Value2 = 0.1245283018867925
$ws.Range("F17").Value2 = 0.01737451737451737
$ws.Range("I17").Value2 = 0.07722007722007722
$ws.Range("J17").Value2 = 0.4691119691119691
$ws.Range("K17").Value2 = 0.09652509652509653
$ws.Range("M17").Value2 = 0.01737451737451737
$ws.Range("O17").Value2 = 0.06177606177606178
$ws.Range("S17").Value2 = 0.1177606177606178
$ws.Range("F18").Value2 = 0.004545454545454545
$ws.Range("H18").Value2 = 0.1454545454545454
$ws.Range("I18").Value2 = 0.08181818181818182
$ws.Range("J18").Value2 = 0.3863636363636364
$ws.Range("K18").Value2 = 0.1181818181818182
$ws.Range("M18").Value2 = 0.01363636363636364
$ws.Range("O18").Value2 = 0.08636363636363636
$ws.Range("S18").Value2 = 0.1636363636363636
$ws.Range("F19").Value2 = 0.01940993788819876
$ws.Range("H19").Value2 = 0.1801242236024845
$ws.Range("I19").Value2 = 0.06521739130434782
$ws.Range("J19").Value2 = 0.3509316770186335
$ws.Range("K19").Value2 = 0.1420807453416149
$ws.Range("M19").Value2 = 0.02329192546583851
$ws.Range("O19").Value2 = 0.06055900621118013
$ws.Range("S19").Value2 = 0.1583850931677019

Write-Output "Applied Auburn_A matrix updates (games pulled March 7)"
